$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 14-19 (ids 10013-10018, French "fra" language entries) ---

$ws.Cells.Item(14, 1).Value = 10013
$ws.Cells.Item(14, 2).Value = "Pré-inscription"
$ws.Cells.Item(14, 3).Value = "Portail Web pour les pré-inscriptions"
$ws.Cells.Item(14, 4).Value = "fra"
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = "superadmin"
$ws.Cells.Item(14, 7).Value = "now()"

$ws.Cells.Item(15, 1).Value = 10014
$ws.Cells.Item(15, 2).Value = "Client dinscription"
$ws.Cells.Item(15, 3).Value = "Application de bureau pour les inscriptions"
$ws.Cells.Item(15, 4).Value = "fra"
$ws.Cells.Item(15, 5).Value = $true
$ws.Cells.Item(15, 6).Value = "superadmin"
$ws.Cells.Item(15, 7).Value = "now()"

$ws.Cells.Item(16, 1).Value = 10015
$ws.Cells.Item(16, 2).Value = "Processeur dinscription"
$ws.Cells.Item(16, 3).Value = "Demande de post-inscription"
$ws.Cells.Item(16, 4).Value = "fra"
$ws.Cells.Item(16, 5).Value = $true
$ws.Cells.Item(16, 6).Value = "superadmin"
$ws.Cells.Item(16, 7).Value = "now()"

$ws.Cells.Item(17, 1).Value = 10016
$ws.Cells.Item(17, 2).Value = "Authentification ID"
$ws.Cells.Item(17, 3).Value = "Application pour lauthentification du fournisseur de services tiers"
$ws.Cells.Item(17, 4).Value = "fra"
$ws.Cells.Item(17, 5).Value = $true
$ws.Cells.Item(17, 6).Value = "superadmin"
$ws.Cells.Item(17, 7).Value = "now()"

$ws.Cells.Item(18, 1).Value = 10017
$ws.Cells.Item(18, 2).Value = "Contrôle didentité"
$ws.Cells.Item(18, 3).Value = "Portail Web pour la configuration dapplications"
$ws.Cells.Item(18, 4).Value = "fra"
$ws.Cells.Item(18, 5).Value = $true
$ws.Cells.Item(18, 6).Value = "superadmin"
$ws.Cells.Item(18, 7).Value = "now()"

$ws.Cells.Item(19, 1).Value = 10018
$ws.Cells.Item(19, 2).Value = "Portail Résident"
$ws.Cells.Item(19, 3).Value = "Portail Web pour les services de génération de post-ID"
$ws.Cells.Item(19, 4).Value = "fra"
$ws.Cells.Item(19, 5).Value = $true
$ws.Cells.Item(19, 6).Value = "superadmin"
$ws.Cells.Item(19, 7).Value = "now()"

# --- Column widths (column A best-fit-ish, column B custom width) ---
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 18.5

# --- View: scroll down a bit and select from row 20 to the bottom of sheet ---
$ws.Range("A20:A1048576").EntireRow.Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10

# --- Page setup: paper size 9 (A4), portrait orientation ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "done"
